# Update Name of Algo
# Apply updated KNN imputation results to column D (header "D") for the
# terrestrial_mammals / combination_3_ABCDF / D / 15 / seed2 scenario.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.382
$ws.Range("D14").Value = -8.077
$ws.Range("D21").Value = -7.9
$ws.Range("D23").Value = -7.505
$ws.Range("D25").Value = -8.388999999999999
